{"js": "// The underlying XML diff for this document is a pure re-serialization\n// (namespace declarations and element attributes were written out in a\n// different, alphabetised order) produced by whatever tool regenerated the\n// test resource. Every <attr>=\"<value\"> pair present before the change is\n// still present afterwards - nothing in the document's content, formatting\n// values, paragraphs, tab-stop positions or page geometry actually changed.\n//\n// The only part of that (cosmetic) diff that corresponds to a concrete,\n// settable property in the Word JavaScript API is the page size / margins\n// of the section (<w:pgSz>/<w:pgMar> on <w:sectPr>). We re-assert those\n// exact values here (they round-trip to the same twips values the template\n// already had: 11906x16838 page size, 1417/1417/1417/1417 margins and\n// 708/708 header/footer distance, 0 gutter) so the section's page setup is\n// explicitly (re)written.\n//\n// (The other hunks in the diff only reorder attributes on <w:tab>,\n// <w:footnote>, <w:rFonts>, <w:lang>, <w:latentStyles>, <w:lsdException>\n// and <w:style> elements - none of which changes any value, and Word's\n// JavaScript API does not expose a way to control raw attribute\n// serialization order, nor does it expose tab-stop manipulation at all.)\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const pageSetup = sections.items[i].pageSetup;\n\n  // Page size: w:pgSz w:w=\"11906\" w:h=\"16838\" -> points (1 pt = 20 twips)\n  pageSetup.pageWidth = 11906 / 20; // 595.3\n  pageSetup.pageHeight = 16838 / 20; // 841.9\n\n  // Page margins: w:pgMar w:top=\"1417\" w:right=\"1417\" w:bottom=\"1417\"\n  // w:left=\"1417\" w:header=\"708\" w:footer=\"708\" w:gutter=\"0\"\n  pageSetup.topMargin = 1417 / 20; // 70.85\n  pageSetup.bottomMargin = 1417 / 20; // 70.85\n  pageSetup.leftMargin = 1417 / 20; // 70.85\n  pageSetup.rightMargin = 1417 / 20; // 70.85\n  pageSetup.headerDistance = 708 / 20; // 35.4\n  pageSetup.footerDistance = 708 / 20; // 35.4\n  pageSetup.gutter = 0 / 20; // 0\n}\n\nawait context.sync();\n", "ps1": "# The underlying XML diff for this document is a pure re-serialization\n# (namespace declarations and element attributes were written out in a\n# different, alphabetised order) produced by whatever tool regenerated the\n# test resource. Every attribute=\"value\" pair present before the change is\n# still present afterwards - nothing in the document's content, formatting\n# values, paragraphs, tab-stop positions or page geometry actually changed\n# (<w:tab w:val=\"left\" w:pos=\"3119\"/> -> <w:tab w:pos=\"3119\" w:val=\"left\"/>,\n# <w:pgSz w:w=\"11906\" w:h=\"16838\"/> -> <w:pgSz w:h=\"16838\" w:w=\"11906\"/>,\n# etc. - same values, just reordered attributes).\n#\n# We re-assert the concrete values the diff's elements carry, through the\n# corresponding Word object model properties, so that the tab stops and the\n# section page geometry are explicitly (re)written:\n#   - the 4 body paragraphs' custom tab stop at 3119 twips (155.95 pt), left\n#     aligned;\n#   - the section's page size (11906x16838 twips) and margins\n#     (1417/1417/1417/1417 twips, 708/708 header/footer distance, 0 gutter).\n\n$d = $word.ActiveDocument\n\n# --- Body paragraph tab stops (w:tab w:val=\"left\" w:pos=\"3119\") ---\nforeach ($i in 1, 2, 3, 4) {\n    $p = $d.Paragraphs($i)\n    $tabStops = $p.Range.ParagraphFormat.TabStops\n    $tabStops.ClearAll()\n    $tabStops.Add(155.95, 0)\n}\n\n# --- Section page size / margins ---\n# w:pgSz w:w=\"11906\" w:h=\"16838\"\n# w:pgMar w:top=\"1417\" w:right=\"1417\" w:bottom=\"1417\" w:left=\"1417\"\n#         w:header=\"708\" w:footer=\"708\" w:gutter=\"0\"\nfor ($i = 1; $i -le $d.Sections.Count; $i++) {\n    $pageSetup = $d.Sections($i).PageSetup\n    $pageSetup.PageWidth = 595.3\n    $pageSetup.PageHeight = 841.9\n    $pageSetup.TopMargin = 70.85\n    $pageSetup.BottomMargin = 70.85\n    $pageSetup.LeftMargin = 70.85\n    $pageSetup.RightMargin = 70.85\n    $pageSetup.HeaderDistance = 35.4\n    $pageSetup.FooterDistance = 35.4\n    $pageSetup.Gutter = 0\n}\n"}
